$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 - new booking SNOW-745596, 2026-02-26
$ws.Range("A7").Value = "SNOW-745596"
$ws.Range("B7").Value = "'2026-02-26"
$ws.Range("C7").Value = "local"
$ws.Range("D7").Value = "loca@a.com"
$ws.Range("E7").Value = "'1212"
$ws.Range("F7").Value = 3
$ws.Range("G7").Value = "Family Ski Package"
$ws.Range("H7").Value = 32000
$ws.Range("I7").Value = 64000
$ws.Range("J7").Value = "Confirmed"
$ws.Range("K7").Value = "'2026-02-23"
$ws.Range("L7").Value = ""

# Row 8 - new booking SNOW-745596, 2026-02-27
$ws.Range("A8").Value = "SNOW-745596"
$ws.Range("B8").Value = "'2026-02-27"
$ws.Range("C8").Value = "local"
$ws.Range("D8").Value = "loca@a.com"
$ws.Range("E8").Value = "'1212"
$ws.Range("F8").Value = 3
$ws.Range("G8").Value = "Family Ski Package"
$ws.Range("H8").Value = 32000
$ws.Range("I8").Value = 64000
$ws.Range("J8").Value = "Confirmed"
$ws.Range("K8").Value = "'2026-02-23"
$ws.Range("L8").Value = ""
